$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.826.04"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "3.386.47"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "558.75"
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("D6").Value = "175.08"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").Value = "3.378.97"
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("E10").Value = "  +8.62%  "
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("D12").Value = "54.06"
$ws.Range("E12").Value = "  +3.56%  "
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +5.12%  "
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").Value = "3.921.18"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "3.377.49"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "64.854.16"
$ws.Range("E19").Value = "  +2.99%  "
$ws.Range("D20").Value = "11.78"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "460.96"
$ws.Range("E22").Value = "  +11.15%  "
$ws.Range("D23").Value = "4.86"
$ws.Range("E23").Value = "  +10.77%  "
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").Value = "86.46"
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "2.94"
$ws.Range("E27").Value = "  +8.95%  "
$ws.Range("D28").Value = "10.83"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").Value = "8.73"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").Value = "30.77"
$ws.Range("E30").Value = "  +6.65%  "
$ws.Range("D31").Value = "6.72"
$ws.Range("D32").Value = "11.46"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "570.23"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "61.17"
$ws.Range("E34").Value = "  +5.86%  "
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("D37").Value = "3.63"
$ws.Range("E37").Value = "  +6.91%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "35.52"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("D42").Value = "3.079.53"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").Value = "0.0416"
$ws.Range("E45").Value = "  +4.52%  "
$ws.Range("E46").Value = "  +5.44%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Value = "3.12"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "137.89"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("D51").Value = "8.23"
$ws.Range("E51").Value = "  +3.40%  "
